$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
        return
    }
    $rng.Text = $new
}

# --- Title ---
Replace-Text "Unraveling the Secrets of Longevity: A Quest for Prolonged Life" "Exploring the Marvels of Science: A Journey Through the Wonders of the Natural World"

# --- Author name ---
Replace-Text "Maria Rodriguez" "Dr. Emily Carson"

# --- Email ---
Replace-Text "mrodriguez@ecobiology" "ecarson@validmail"
Replace-Text "edu" "net"

# --- Body paragraph: intro/science block ---
Replace-Text "The human journey has long been adorned with stories of eternal youth and elusive immortality" "Science, the systematic study of the natural world, unveils the intricate tapestry of life and the enigmatic symphony of the universe"

Replace-Text " From ancient alchemists experimenting with elixirs of life to tales of ageless gods and revitalizing waters, our fascination with prolonged existence echoes in myths and legends" " It empowers us to unravel the mysteries of existence, unmask the enigmatic secrets of nature's dance, and decipher the quantum enigma of reality"

Replace-Text " In modern times, scientific advancements have transformed the study of aging, unveiling avenues to potentially extend the human lifespan. Embarking on this quest for longevity, researchers have ventured into understanding the complex interplay of genes, nutrition, cellular mechanisms, and environmental factors that influence the aging process" " Science's relentless pursuit of knowledge expands our understanding, enriches our lives, and propels humanity forward"

# --- Body paragraph: chemistry block ---
Replace-Text "Our understanding of aging has evolved dramatically over the years" "In the realm of chemistry, we delve into the intricate dance of atoms and molecules, unraveling the secrets of chemical reactions and compounds"

Replace-Text " Scientists have identified genetic factors that play a vital role in longevity, unraveling the intricacies of inherited biological mechanisms" " We uncover the fundamental principles governing the composition of matter, the interactions between substances, and the remarkable diversity of chemical phenomena"

Replace-Text " Nutrition has also emerged as a key determinant of lifespan, highlighting the importance of balanced diets and the potential impact of specific nutrients in extending cellular health. Delving into the cellular realm, research has revealed intricate pathways and molecular processes that contribute to aging, setting the stage for breakthroughs that could potentially slow or reverse these processes" " Chemistry unravels the mesmerizing symphony of the molecular world, revealing the incredible complexity and beauty inherent in the natural world"

# --- Body paragraph: biology block ---
Replace-Text "Environmental influences on aging cannot be overlooked" "Biology, the study of life, unveils the symphony of life's intricacies, from the microscopic realm of cells to the vast array of organisms that inhabit our planet"

Replace-Text " Exposure to pollutants, stress, and various lifestyle factors such as smoking and obesity have a profound effect on the aging trajectory" " We explore the intricate mechanisms underlying life processes, unravel the mysteries of DNA and genetic inheritance, and decipher the remarkable adaptation of living organisms to their diverse environments"

Replace-Text " Understanding the intricate web of interactions between the internal biology of organisms and the external environment holds the key to optimizing health and longevity" " Biology unveils the profound interconnectedness of life, revealing the symphony of interactions within ecosystems and the delicate balance of nature"

# --- Summary paragraph ---
Replace-Text "This essay delved into the fascinating journey of exploring longevity, highlighting the contributions of genetic, nutritional, cellular, and environmental factors in determining lifespan" "Through the study of science, we embark on a journey of exploration and enlightenment"

Replace-Text " We examined the intricate dance between our genes, the food we consume, the inner workings of our cells, and the environment we inhabit" " Chemistry unveils the intricate dance of atoms and molecules, unmasking the secrets of matter's diversity"

Replace-Text " The quest for longevity continues to captivate the human imagination, inspiring scientific inquiry and fueling our desire to transcend the constraints of mortality" " Biology unravels the symphony of life, revealing the interconnectedness and complexity of living organisms. These fields of study empower us to comprehend the wonders of the natural world, address global challenges, improve human health, and create innovative technologies. As we continue to delve into the enigmatic tapestry of science, we unlock the mysteries of the universe and create a brighter future for humanity"

# --- Add a new empty paragraph at the end of the document body ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endRng = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$endRng.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')
